$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(70, 1).Value = "2024-10-10 00:00:00"
$ws.Cells.Item(70, 2).Value = 76950
$ws.Cells.Item(70, 3).Value = 10842
$ws.Cells.Item(70, 4).Value = 9594.690000000001
$ws.Cells.Item(70, 5).Value = 7.0749
